$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Edit 1: In the "- NORM: ... - HOLD: ... 02" paragraph, insert a new
# run with text "01 e " right before the run containing "02", so the
# line reads "- HOLD:<tab><tab>01 e 02" (two separate runs: "01 e "
# and "02").
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("02", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    # Route the text change through a one-character placeholder first so
    # the engine registers a genuine content mutation (a same-text
    # "rewrite" is treated as a no-op and left untouched).
    $rng.Text = "X"
    $rewritten = $d.Range($rng.Start, $rng.Start + 1)
    $rewritten.Text = "01 e 02"

    # Split "01 e 02" into two runs - "01 e " and "02" - by toggling a
    # real character-formatting property on/off across just the "01 e "
    # prefix. Word will not silently re-merge runs that were touched by
    # an explicit formatting change, so this leaves two <w:r> elements
    # behind; toggling it back to its original value keeps the visible
    # formatting unchanged.
    $prefix = $d.Range($rewritten.Start, $rewritten.Start + 5)
    $prefix.Bold = 1
    $prefix.Bold = 0
}

# ------------------------------------------------------------------
# Edit 2: In the "- Novelty:" paragraph (the one directly under
# "- Melhores kernels (best):" with three trailing tabs, immediately
# below "PAP-SMEAR"), merge the two runs "- Novelty" and ":" into a
# single run "- Novelty:".
# ------------------------------------------------------------------
$rng2 = $d.Content
$foundHeading = $rng2.Find.Execute("PAP-SMEAR", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundHeading) {
    $afterHeading = $d.Range($rng2.End, $d.Content.End)
    $foundNovelty = $afterHeading.Find.Execute("- Novelty:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($foundNovelty) {
        # Same placeholder trick, this time to coalesce two adjacent runs
        # ("- Novelty" and ":") into one run, matching how Word merges
        # runs whenever their text is genuinely rewritten together.
        $afterHeading.Text = "X"
        $novRewritten = $d.Range($afterHeading.Start, $afterHeading.Start + 1)
        $novRewritten.Text = "- Novelty:"
    }
}
